$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the "This spike serves..." paragraph into a single run.
#    (Removes the spell-check proofErr wrappers around "on-premise" /
#    "serverless" left over from Word's spell checker.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This spike serves to evaluate the performance of and difficulty developing an on-premise application. The results will be used to compare against a serverless implementation.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This spike serves to evaluate the performance of and difficulty developing an on-premise application. The results will be used to compare against a serverless implementation.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Collapse the "A node.js based program..." paragraph into a single run.
#    (Removes the spell-check proofErr wrappers around "json".)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "A node.js based program that reads the data from the Arduino using a motion detector and pushes it to the firebase. The data that the program pushes should be in the json format with 3 attributes:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "A node.js based program that reads the data from the Arduino using a motion detector and pushes it to the firebase. The data that the program pushes should be in the json format with 3 attributes:",
    2) | Out-Null

# The old "_GoBack" bookmark (which sat right after "...serverless
# implementation.") is dropped here; a fresh one is added at the new edit
# location created below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) Add a new paragraph after "Planning notes:  " describing how the work
#    was divided between the team members.
# ---------------------------------------------------------------------------
$part1 = "To achieve this spike we have decided to have one team member create the first component that will read data from the board, and have the other implement "
$part2 = "the on premise node.js function, component2, that will use nodemailer to send the emails from the node app"
$part3 = ". This will allow one us to finish the spike at the same time, and both work together to complete the report.  "
$fullText = $part1 + $part2 + $part3

$d.Content.Find.Execute(
    "Planning notes:  ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Planning notes:  ^p" + $fullText,
    2) | Out-Null

$newPara = $d.Paragraphs.Last
$paraStart = $newPara.Range.Start

# Split "part1" from "part2" into separate runs using a throwaway bookmark
# (runs with identical formatting would otherwise be re-combined); the
# bookmark itself is removed again immediately afterwards.
$splitPos = $paraStart + $part1.Length
$d.Bookmarks.Add("_TempRunSplit", $d.Range($splitPos, $splitPos))

# Re-create "_GoBack" between "part2" and "part3", matching the new edit
# location (this also forces that run boundary to persist).
$goBackPos = $paraStart + $part1.Length + $part2.Length
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))

if ($d.Bookmarks.Exists("_TempRunSplit")) {
    $d.Bookmarks("_TempRunSplit").Delete()
}

# Match the body-text paragraph style used elsewhere in the document.
$newParaFmt = $d.Paragraphs.Last.Range.ParagraphFormat
$newParaFmt.SpaceAfter = 0.7
$newParaFmt.LeftIndent = -0.25
